$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# D-column values are forced to Text format first so Excel does not
# reinterpret numeric-looking strings (stripping trailing zeros, etc.).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.415.70"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.369.41"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.62"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.85"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.72"
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.46"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.980"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.730.58"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.31"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.369.85"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.421.04"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.10"
$ws.Range("E19").Value = "  +7.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000106"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.22"
$ws.Range("E21").Value = "  -5.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.25"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.52"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.82"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.13"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").Value = "  -5.57%  "
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0980"
$ws.Range("E30").Value = "  +5.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.30"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.12"
$ws.Range("E32").Value = "  -4.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "166.41"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.71"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.97"
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").Value = "  +8.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0354"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.83"
$ws.Range("E42").Value = "  -5.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.90"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("E44").Value = "  -4.98%  "
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.65"
$ws.Range("E46").Value = "  -8.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.815.36"
$ws.Range("E47").Value = "  +9.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.02"
$ws.Range("E48").Value = "  +5.43%  "
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.95"
$ws.Range("E51").Value = "  -6.57%  "
